# Apply the "Correzione di alcune exit condition" edits.
# Uses Find/Execute (wdReplaceOne) against $d.Content so each call re-scans
# the whole story range; every search string is unique in the document.

$d = $word.ActiveDocument
$apos = [char]0x2019

function Replace-Text($find, $replace) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $find"
    }
}

# 1) Use case name cell
$find1 = "Visualizzazione di un modello di un auto"
$repl1 = "Visualizzazione di un modello scelto dall" + $apos + "utente"
Replace-Text $find1 $repl1

# 2) Descrizione cell
$find2 = "visualizzare una pagina che contiene le caratteristiche del modello di un" + $apos + " auto scelta dall" + $apos + "utente"
$repl2 = "visualizzare una pagina che contiene le caratteristiche del modello scelto dall" + $apos + "utente"
Replace-Text $find2 $repl2

# 3) Entry condition cell
$find3 = "clicca su un modello per scoprire le sue caratteristiche."
$repl3 = "clicca sull" + $apos + "anteprima di un modello per visualizzare una pagina che contiene le caratteristiche del modello."
Replace-Text $find3 $repl3

# 4) Exit condition (On success) cell
$find4 = "L" + $apos + "utente visualizza una pagina che contiene le caratteristiche del modello dell" + $apos + "auto scelta."
$repl4 = "L" + $apos + "utente visualizza una pagina che contiene le caratteristiche del modello scelto."
Replace-Text $find4 $repl4

# 5) Exit condition (On failure) cell
$find5 = "L" + $apos + "utente non riesce a comunicare col server e non riesce a vedere la pagina."
$repl5 = "L" + $apos + "utente non riesce a vedere la pagina."
Replace-Text $find5 $repl5

# 6) Scenario step cell
$find6 = "al sistema di mostrargli le caratteristiche di un modello di un" + $apos + "auto da lui scelto facendo una richiesta http al server di sistema."
$repl6 = "al sistema di mostrargli le caratteristiche di un modello da lui scelto facendo una richiesta http al server di sistema."
Replace-Text $find6 $repl6

Write-Output "done"
